$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NOV-2020")

# Copy formatting (styles/borders/number formats) from the last existing
# data row (62) down into the new rows we are about to populate (63-65),
# then onto the trailing blank rows (66-69), matching the template used
# throughout the sheet.
$ws.Range("A62:G62").Copy() | Out-Null
$ws.Range("A63:G65").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(63).RowHeight = $ws.Rows.Item(62).RowHeight

# Row 63: 29 | 29-Nov-2020 | RPA RLOGIC | (reuse existing task text) | 90% | WIP | supported on Sunday
$ws.Cells.Item(63, 1).Value = 29
$ws.Cells.Item(63, 2).Value = 44164
$ws.Cells.Item(63, 3).Value = "RPA RLOGIC"
$ws.Cells.Item(63, 4).Value = $ws.Cells.Item(62, 4).Value
$ws.Cells.Item(63, 5).Value = 0.9
$ws.Cells.Item(63, 6).Value = "WIP"
$ws.Cells.Item(63, 7).Value = "supported on Sunday"

# Row 64: 30 | 30-Nov-2020 | RPA RLOGIC | Suppoted to RLOGIC for P&L | 100% | Completed |
$ws.Cells.Item(64, 1).Value = 30
$ws.Cells.Item(64, 2).Value = 44165
$ws.Cells.Item(64, 3).Value = "RPA RLOGIC"
$ws.Cells.Item(64, 4).Value = "Suppoted to RLOGIC for P&L"
$ws.Cells.Item(64, 5).Value = 1
$ws.Cells.Item(64, 6).Value = "Completed"

# Row 65: 31 | 30-Nov-2020 | RPA GSPN | Corrections received... | 100% | Completed |
$ws.Cells.Item(65, 1).Value = 31
$ws.Cells.Item(65, 2).Value = 44165
$ws.Cells.Item(65, 3).Value = "RPA GSPN"
$ws.Cells.Item(65, 4).Value = "Corrections received for Warranty task and it has been completed from download to upload process"
$ws.Cells.Item(65, 5).Value = 1
$ws.Cells.Item(65, 6).Value = "Completed"

# Trailing blank rows 66-69 reuse the plain (borderless) look applied
# elsewhere in the sheet to separator rows.
$ws.Range("A66:G69").Style = "Normal"
$ws.Range("A66:G69").HorizontalAlignment = -4131
$ws.Range("A66:G69").VerticalAlignment = -4108

# Row 70 closes the table out with a thin bottom border.
$ws.Range("A70:G70").Borders.Item(9).LineStyle = 1
$ws.Range("A70:G70").Borders.Item(9).Weight = 2
$ws.Range("A70:G70").HorizontalAlignment = -4131
$ws.Range("A70:G70").VerticalAlignment = -4108

# Widen column G to fit the newly added comment text.
$ws.Columns.Item(7).ColumnWidth = 18.21875

$ws.Application.GoTo($ws.Range("D68"))
